$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value of 45186 (2023-09-17)
# for rows 2 through 72; update it to 45188 (2023-09-19) for all of them.
$ws.Range("C2:C72").Value = 45188
